# Reorder the names listed in the "Recorded By" (column G) cells.
# - "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# - "System, backup@backdoor.com, system" -> "System, system, backup@backdoor.com"
# Only cells whose text matches exactly one of these two source strings are touched;
# every other value in the column (e.g. already-reordered rows, or rows with other
# recorder combinations) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
